# Fix Q calculation: add two missing generator rows to the "gen" sheet,
# and correct the from_bus/to_bus/length_km values for two lines on the
# "line" sheet (Line6 and Line7) that fed into the wrong Q numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "gen" sheet: append rows 3 and 4 (new generators), matching the
# formatting (bold/bordered index column) already used by row 2.
# ---------------------------------------------------------------------
$genWs = $wb.Worksheets.Item("gen")

# Row 3
$genWs.Range("A2").Copy($genWs.Range("A3"))
$genWs.Range("A3").Value = 1
$genWs.Range("C3").Value = 8
$genWs.Range("D3").Value = 17.25
$genWs.Range("E3").Value = 1
$genWs.Range("I3").Value = 1
$genWs.Range("J3").Value = $false
$genWs.Range("K3").Value = $true
$genWs.Range("L3").Value = 0

# Row 4
$genWs.Range("A2").Copy($genWs.Range("A4"))
$genWs.Range("A4").Value = 2
$genWs.Range("C4").Value = 7
$genWs.Range("D4").Value = 42
$genWs.Range("E4").Value = 1
$genWs.Range("I4").Value = 1
$genWs.Range("J4").Value = $false
$genWs.Range("K4").Value = $true
$genWs.Range("L4").Value = 0

# ---------------------------------------------------------------------
# "line" sheet: fix from_bus / to_bus / length_km for Line6 (row 8)
# and Line7 (row 9).
# ---------------------------------------------------------------------
$lineWs = $wb.Worksheets.Item("line")

# Line6 (row 8)
$lineWs.Range("D8").Value = 8
$lineWs.Range("E8").Value = 2
$lineWs.Range("F8").Value = 15.12

# Line7 (row 9)
$lineWs.Range("D9").Value = 7
$lineWs.Range("E9").Value = 4
$lineWs.Range("F9").Value = 18.33
